$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RetanqueoDigiCredito")
$ws.Activate()

# Update the DigiCredito sample-data row (row 2) with the new test values.
# Order matters: it reproduces the order new shared strings were authored in.
$ws.Range("V2").Value = '"86292"'
$ws.Range("M2").Value = '"OSCAR"'
$ws.Range("A2").Value = '"7500000"'
$ws.Range("B2").Value = '"10092369"'
$ws.Range("C2").Value = '"68003"'
$ws.Range("F2").Value = '"50"'
$ws.Range("G2").Value = '"50"'
$ws.Range("H2").Value = '"8700000"'
$ws.Range("I2").Value = '"250000"'
$ws.Range("J2").Value = '"300000"'

# Move the selection/active cell (no more scrolled to topLeftCell G1;
# the cursor now rests on L4 instead of W4).
$ws.Range("L4").Select()
